$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.660.25'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.491'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0620'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0834'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.862.80'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.630.78'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.639.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.03'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.04'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.37%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.32'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.12'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.42'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.43'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.68'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.52%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.165.10'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.795'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.771.61'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.38'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.56'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0105'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.61'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.56'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.09%  '
